$d = $word.ActiveDocument

function Frag($inner) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $inner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Target range: from the start of paragraph 2 (the old 'Vou falar...' paragraph)
# through the end of the document's last paragraph -- this covers everything
# that changes, leaving paragraph 1 ('Seguranca da informacao') untouched.
$startPara = $d.Paragraphs(2)
$totalBefore = $d.Paragraphs.Count
$endPara = $d.Paragraphs($totalBefore)
$r = $d.Range($startPara.Range.Start, $endPara.Range.End)

$batch0 = '<w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>1º SLIDE</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Vou falar um bocado da parte conceitual </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">e vou dar </w:t>
      </w:r>
      <w:r>
        <w:t>dois exemplos de como é que se consegue garantir a segurança da informação</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>º SLIDE</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">A informação é um conjunto de dados que inseridos num contexto passam a ter significado para pessoas ou organizações. Isto pode ser, informações financeiras, segredos políticos, arquivos confidenciais, dados pessoais, entre tantas outras coisas. Para proteger estas informações surge a segurança da informação. </w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t>Onde existe informação?  - Existe informação na comunicação, troca de mensagens, emails, conversa entre pessoas, etc.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>3</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>º SLIDE</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Para garantir </w:t>
      </w:r>
      <w:r>
        <w:t>segurança na troca destas mensagens, existem 3 propriedades mais importantes na segurança da informação.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t>Confidencialidade – Somente pessoas autorizadas têm acesso à mensagem/informação.</w:t>
      </w:r>
    </w:p>'
$null = $r.InsertXML((Frag $batch0))
$insertedCount = 8

$batch1 = '<w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t>Integridade – Garantia de que a mensagem não foi alterada ou modificada durante a comunicação. O recetor recebe a mensagem do emissor sem alterações.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:t>Disponibilidade – Garantia de que a informação estará pronta para usar por pessoas autorizadas quando for necessária.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>4</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>º</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t xml:space="preserve"> SLIDE</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Há pessoas que dizem que sim, outras o contrário. Pela minha pesquisa o que encontrei, foi que não existe 100% de segurança porque para que isso existisse teria de haver um grande investimento por parte das empresas. Desta forma é preciso balancear entre o risco da informação que possa ser exposta e o preço que estamos dispostos a pagar para garantir que a informação </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">não </w:t>
      </w:r>
      <w:r>
        <w:t>seja revelada para quem nós não queremos.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Se a informação não valer muito, não compensa ter segurança da informação 100% pois seria muito caro. Portanto é preciso ter isso em conta.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>5º SLIDE</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Como podemos garantir que a aplicação móvel ou outra aplicação que interaja com a casa não esteja comprometida a fim de um atacante poder saber se a pessoa está ou não em casa com a localização do telemóvel e a localização da casa, conseguir </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">saber se as portas ou janelas estão </w:t>
      </w:r>
      <w:r>
        <w:t>tranca</w:t>
      </w:r>
      <w:r>
        <w:t>das</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> e destranca</w:t>
      </w:r>
      <w:r>
        <w:t>das</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, controlar a temperatura, as luzes, os eletrodomésticos que possamos ter na cozinha, </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">saber se o alarme está </w:t>
      </w:r>
      <w:r>
        <w:t>liga</w:t>
      </w:r>
      <w:r>
        <w:t>do</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> ou desliga</w:t>
      </w:r>
      <w:r>
        <w:t>do</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, enfim uma panóplia de ataques possíveis. Para evitar isto tudo, tem de haver um grande sistema de segurança da informação para não haver vulnerabilidades </w:t>
      </w:r>
      <w:r>
        <w:t>nos diversos dispositivos que estejam conectados à internet e que a partir dele um atacante possa tirar partido.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>6º SLIDE</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Como podemos garantir que uma cidade inteligente não exponha informações privadas dos seus moradores? </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Ou prevenir que um atacante tenha acesso às informações dos semáforos de </w:t>
      </w:r>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>trânsito, ou às câmaras de vigilância dessa cidade? Para evitar isto também é necessário haver um sistema de segurança de informação robusto e capaz de não ter vulnerabilidades que possam comprometer a vida dos habitantes dessa cidade.</w:t>
      </w:r>
    </w:p>'
$lastInsertedPara = $d.Paragraphs(1 + $insertedCount)
$rNext = $d.Range($lastInsertedPara.Range.End, $lastInsertedPara.Range.End)
$null = $rNext.InsertXML((Frag $batch1))
$insertedCount = $insertedCount + 8

$batch2 = '<w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>7º SLIDE</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t>Portanto, se voltarmos a olhar para o exemplo 1 e compararmos com o exemplo 2 em termos de como garantir a segurança da informação</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> conseguimos perceber que não existem quaisquer diferenças. O problema é sempre o mesmo. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">A informação que os atacantes possam obter é que é diferente. E esta pode ser mais ou menos valiosa. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Não vale a pena estar a discutir a segurança da informação nas casas inteligentes ou nas cidades inteligentes, ou onde quer que seja. O problema é sempre o mesmo. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">É preciso também ter em conta a balança demonstrada anteriormente. O risco de expor a informação versus o preço a pagar para não correr esse risco. Portanto, </w:t>
      </w:r>
      <w:r>
        <w:t>para</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> pode</w:t>
      </w:r>
      <w:r>
        <w:t>r</w:t>
      </w:r>
      <w:r>
        <w:t>mos garantir que a informação é assegurada da melhor maneira de modo a ser privada para quem tem de ser</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> é preciso, mais uma vez, ter em conta a balança e uma grande equipa de engenheiros de segurança da informação</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>'
$lastInsertedPara = $d.Paragraphs(1 + $insertedCount)
$rNext = $d.Range($lastInsertedPara.Range.End, $lastInsertedPara.Range.End)
$null = $rNext.InsertXML((Frag $batch2))
$insertedCount = $insertedCount + 2

Write-Host "Final paragraph count:" $d.Paragraphs.Count
